# Generate Report for Handback
# Updates the handback-status report with freshly generated handoff/handback
# timestamps for the 443ec728-68da-42c3-ad74-50858ff18bbf file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for the en-US file row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 12:59:31"

# --- zh-cn sheet: refresh handoff / handback datetimes for the same file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 12:59:26"
$wsZhCn.Range("K2").Value = "2016-09-05 12:59:43"

# --- de-de sheet: refresh handoff / handback datetimes for the same file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 12:59:31"
$wsDeDe.Range("K2").Value = "2016-09-05 12:59:52"
